$d = $word.ActiveDocument

# Simple one-for-one text replacements within the worksheet table cells.
# Each "old" value is unique in the document, so a plain Find/Replace
# targets exactly the intended cell.
$replacements = @(
    @("37÷6=6, 1", "45÷9=5, 0"),
    @("24÷6=4, 0", "52÷5=10, 2"),
    @("41÷2=20, 1", "26÷6=4, 2"),
    @("42÷6=7, 0", "13÷9=1, 4"),
    @("68÷5=13, 3", "33÷8=4, 1"),
    @("18÷6=3, 0", "79÷7=11, 2"),
    @("56÷2=28, 0", "45÷6=7, 3"),
    @("28÷7=4, 0", "88÷9=9, 7"),
    @("14÷6=2, 2", "57÷9=6, 3"),
    @("21÷8=2, 5", "61÷3=20, 1"),
    @("21÷5=4, 1", "43÷8=5, 3"),
    @("36÷7=5, 1", "52÷5=10, 2"),
    @("11÷3=3, 2", "22÷4=5, 2"),
    @("33÷9=3, 6", "19÷7=2, 5"),
    @("97÷2=48, 1", "91÷5=18, 1"),
    @("96÷8=12, 0", "70÷6=11, 4"),
    @("88÷5=17, 3", "16÷7=2, 2"),
    @("97÷4=24, 1", "35÷4=8, 3"),
    @("12÷8=1, 4", "21÷9=2, 3"),
    @("82÷9=9, 1", "61÷7=8, 5"),
    @("72÷9=8, 0", "88÷6=14, 4"),
    @("41÷5=8, 1", "67÷7=9, 4"),
    @("52÷9=5, 7", "59÷2=29, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

# The row that used to read:
#   96÷8=12, 0 | 88÷5=17, 3 | 97÷4=24, 1 | 82÷3=27, 1 | 60÷2=30, 0
# now reads:
#   70÷6=11, 4 | 16÷7=2, 2 | 35÷4=8, 3 | 34÷5=6, 4 | 82÷3=27, 1
# i.e. a new "34÷5=6, 4" cell is inserted before the trailing "82÷3=27, 1"
# cell, and the old trailing "60÷2=30, 0" cell is dropped -- net effect:
# the 4th and 5th cells of that (already-updated) row change in place.
$t = $d.Tables.Item(1)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $row = $t.Rows.Item($r)
    if ($row.Cells.Count -ge 5 -and $row.Cells.Item(4).Range.Text -like "82÷3=27, 1*" -and $row.Cells.Item(5).Range.Text -like "60÷2=30, 0*") {
        $row.Cells.Item(4).Range.Text = "34÷5=6, 4"
        $row.Cells.Item(5).Range.Text = "82÷3=27, 1"
    }
}
